# ------------------------------------------------------------------------
# Catch_Trust_60.xlsx - "Update of all scripts and data"
#
# The species catch list is kept alphabetically sorted within each
# Survey/Area/Station/Gear group. This update:
#   * adds a missing "Arnoglossus thori" catch record to the 1-RAP group
#     (rows 2-20), which shifts the remaining alphabetical rows down by one
#   * drops the erroneous "Scyliorhinus canicula" record from that same
#     group (it lands on the last slot of the shifted block)
#   * drops the duplicate "Arnoglossus thori" record that used to live in
#     the big 2-RAP "other species" group (rows 33-72), shifting that
#     block up by one
#   * refreshes the raising-factor (RF, column I) for every row of that
#     2-RAP "other species" group
#   * corrects several Numb (column H) counts in that same group
#
# Net effect: the sheet shrinks from 72 to 71 rows (header + 70 records).
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> { column letter -> new value } for every row whose contents change.
$rowUpdates = [ordered]@{
    4 = @{ 'E' = 'Arnoglossus thori'; 'F' = 'ARNOTHO'; 'G' = 0.003; 'H' = 1 }
    5 = @{ 'E' = 'Blennius ocellaris'; 'F' = 'BLENOCE'; 'G' = 0.104; 'H' = 4 }
    6 = @{ 'E' = 'Buglossidium luteum'; 'F' = 'BUGLLUT'; 'G' = 0.037 }
    7 = @{ 'E' = 'Citharus linguatula'; 'F' = 'CITHMAC'; 'G' = 0.076 }
    8 = @{ 'E' = 'Eutrigla gurnardus'; 'F' = 'EUTRGUR'; 'G' = 0.046; 'H' = 1 }
    9 = @{ 'E' = 'Lepidotrigla cavillone'; 'F' = 'LEPTCAV'; 'G' = 0.023; 'H' = 3 }
    10 = @{ 'E' = 'Microchirus ocellatus'; 'F' = 'MICUOCE'; 'G' = 0.038; 'H' = 1 }
    11 = @{ 'E' = 'Mullus barbatus'; 'F' = 'MULLBAR'; 'G' = 0.076; 'H' = 2 }
    12 = @{ 'E' = 'Pagellus erythrinus'; 'F' = 'PAGEERY'; 'G' = 0.345; 'H' = 16 }
    13 = @{ 'E' = 'Pecten jacobaeus'; 'F' = 'PECTJAC'; 'G' = 0.795; 'H' = 14 }
    14 = @{ 'E' = 'Raja clavata'; 'F' = 'RAJACLA'; 'G' = 3.42; 'H' = 1 }
    15 = @{ 'E' = 'Scorpaena notata'; 'F' = 'SCORNOT'; 'G' = 0.896; 'H' = 23 }
    16 = @{ 'E' = 'Scorpaena scrofa'; 'F' = 'SCORSCO'; 'G' = 0.216; 'H' = 2 }
    33 = @{ 'I' = 127.3592682926829 }
    34 = @{ 'I' = 127.3592682926829 }
    35 = @{ 'I' = 127.3592682926829 }
    36 = @{ 'I' = 127.3592682926829 }
    37 = @{ 'I' = 127.3592682926829 }
    38 = @{ 'E' = 'Ascidia mentula'; 'F' = 'ASCIMEN'; 'G' = 0.002; 'H' = 1; 'I' = 127.3592682926829 }
    39 = @{ 'E' = 'Ascidia sp'; 'F' = 'ASCIDSP'; 'G' = 0.003; 'I' = 127.3592682926829 }
    40 = @{ 'E' = 'Ascidiacea nd'; 'F' = 'ASCIDND'; 'G' = 0.22; 'H' = -1; 'I' = 127.3592682926829 }
    41 = @{ 'E' = 'Biological discard'; 'F' = 'BIOLDIS'; 'G' = 1.729; 'H' = -1; 'I' = 127.3592682926829 }
    42 = @{ 'E' = 'Bryozoa nd'; 'F' = 'BRYOZND'; 'G' = 0.003; 'H' = 6; 'I' = 127.3592682926829 }
    43 = @{ 'E' = 'Calyptraea chinensis'; 'F' = 'CALICHI'; 'G' = 0.001; 'H' = 1; 'I' = 127.3592682926829 }
    44 = @{ 'E' = 'Chlamys varia'; 'F' = 'CHLAVAR'; 'G' = 0.008; 'H' = 4; 'I' = 127.3592682926829 }
    45 = @{ 'E' = 'Coralligenous concretions'; 'F' = 'CORACON'; 'G' = 0.189; 'H' = -1; 'I' = 127.3592682926829 }
    46 = @{ 'E' = 'Dromia personata'; 'F' = 'DROMPER'; 'G' = 0.008; 'H' = 3; 'I' = 127.3592682926829 }
    47 = @{ 'E' = 'Eggs of Scyliorhinus'; 'F' = 'EGGSSCY'; 'G' = 0.001; 'H' = 1; 'I' = 127.3592682926829 }
    48 = @{ 'E' = 'Eurynome aspera'; 'F' = 'EURYASP'; 'G' = 0.006; 'H' = 7; 'I' = 127.3592682926829 }
    49 = @{ 'E' = 'Galathea sp'; 'F' = 'GALATSP'; 'G' = 0.001; 'H' = 1; 'I' = 127.3592682926829 }
    50 = @{ 'E' = 'Holothuria forskali'; 'F' = 'HOLOFOR'; 'G' = 0.551; 'H' = 8; 'I' = 127.3592682926829 }
    51 = @{ 'E' = 'Holothuria tubulosa'; 'F' = 'HOLOTUB'; 'G' = 0.066; 'H' = 1; 'I' = 127.3592682926829 }
    52 = @{ 'E' = 'Inachus dorsettensis'; 'F' = 'INACDOR'; 'G' = 0.005; 'H' = 4; 'I' = 127.3592682926829 }
    53 = @{ 'E' = 'Laevicardium oblongum'; 'F' = 'LAEVCAR'; 'G' = 0.001; 'H' = 2; 'I' = 127.3592682926829 }
    54 = @{ 'E' = 'Macropodia rostrata'; 'F' = 'MACRROS'; 'H' = 1; 'I' = 127.3592682926829 }
    55 = @{ 'E' = 'Marthasterias glacialis'; 'F' = 'MARTGLA'; 'G' = 0.074; 'I' = 127.3592682926829 }
    56 = @{ 'E' = 'Microcosmus sulcatus'; 'F' = 'MICSULC'; 'G' = 0.285; 'H' = 90; 'I' = 127.3592682926829 }
    57 = @{ 'E' = 'Nassarius nitidus'; 'F' = 'NASSNIT'; 'G' = 0.001; 'H' = 1; 'I' = 127.3592682926829 }
    58 = @{ 'E' = 'Ophiothrix fragilis'; 'F' = 'OPHIFRA'; 'G' = 0.344; 'H' = 430; 'I' = 127.3592682926829 }
    59 = @{ 'E' = 'Ophiotrix quinquemaculata'; 'F' = 'OPHIQUI'; 'G' = 0.152; 'H' = 177; 'I' = 127.3592682926829 }
    60 = @{ 'E' = 'Ophiura albida'; 'F' = 'OPHIALB'; 'G' = 0.002; 'H' = 12; 'I' = 127.3592682926829 }
    61 = @{ 'E' = 'Ophiura ophiura'; 'F' = 'OPHIOPH'; 'G' = 0.007; 'H' = 15; 'I' = 127.3592682926829 }
    62 = @{ 'E' = 'Pagurus cuanensis'; 'F' = 'PAGUCUA'; 'G' = 0.001; 'H' = 1; 'I' = 127.3592682926829 }
    63 = @{ 'E' = 'Pagurus prideauxi'; 'F' = 'PAGUPRI'; 'G' = 0.006; 'H' = 4; 'I' = 127.3592682926829 }
    64 = @{ 'E' = 'Parthenope massena'; 'F' = 'PARTMAS'; 'G' = 0.003; 'I' = 127.3592682926829 }
    65 = @{ 'E' = 'Pilumnus spinifer'; 'F' = 'PILUSPI'; 'G' = 0.002; 'H' = 3; 'I' = 127.3592682926829 }
    66 = @{ 'E' = 'Porifera nd'; 'F' = 'PORIFND'; 'G' = 0.796; 'H' = -1; 'I' = 127.3592682926829 }
    67 = @{ 'E' = 'Psammechinus microtuberculatus'; 'F' = 'PSAMMIC'; 'G' = 0.026; 'H' = 13; 'I' = 127.3592682926829 }
    68 = @{ 'E' = 'Pteria hirundo'; 'F' = 'PTERHIR'; 'G' = 0.006; 'H' = 1; 'I' = 127.3592682926829 }
    69 = @{ 'E' = 'Shells NA'; 'F' = 'SHELLS'; 'G' = 0.59; 'H' = -1; 'I' = 127.3592682926829 }
    70 = @{ 'E' = 'Wood NA'; 'F' = 'WOOD'; 'G' = 0.15; 'H' = -1; 'I' = 127.3592682926829 }
    71 = @{ 'E' = 'Xanto pilipes'; 'F' = 'XANTPHI'; 'G' = 0.001; 'H' = 2; 'I' = 127.3592682926829 }
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cellValues = $rowUpdates[$rowNum]
    foreach ($col in $cellValues.Keys) {
        $colIndex = [int][char]$col - [int][char]"A" + 1
        $ws.Cells.Item($rowNum, $colIndex).Value2 = $cellValues[$col]
    }
}

# The whole data block shrank by one record overall (one new row added to the
# 1-RAP group, two rows removed from the 2-RAP "other species" group); remove
# the now-obsolete trailing row so the used range/dimension shrinks to K71.
$ws.Rows.Item(72).Delete()
